# Updated renewables and battery capital and operating costs to 2024 ATB and 2022 USD
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row labels: "2020 ..." -> "2022 ..." ---
$ws.Range("A9").Value  = "2022 CapEx"
$ws.Range("A15").Value = "2022 OpEx ($/kw-yr)"
$ws.Range("A22").Value = "2022 PV base installed cost"
$ws.Range("A28").Value = "2022 PV OpEx"

$cols = @("B", "C", "D", "E", "F")

# row number -> new values for columns B..F
$rowData = @{
    9  = @(1666, 1666, 1803, 2335, 1666)   # Wind CapEx (2022)
    10 = @(1569, 1569, 1703, 2151, 1569)   # Wind CapEx (2025)
    11 = @(1408, 1408, 1537, 1844, 1408)   # Wind CapEx (2030)
    12 = @(1335, 1335, 1457, 1749, 1335)   # Wind CapEx (2035)
    13 = @(1115, 1115, 1219, 1463, 1115)   # Wind CapEx (2050)

    15 = @(32, 32, 30, 32, 32)             # Wind OpEx (2022)
    16 = @(31, 31, 29, 31, 31)             # Wind OpEx (2025)
    17 = @(29, 29, 27, 29, 29)             # Wind OpEx (2030)
    18 = @(28, 28, 26, 28, 28)             # Wind OpEx (2035)
    19 = @(25, 25, 23, 25, 25)             # Wind OpEx (2050)

    22 = @(1483, 1483, 1483, 1483, 1483)   # PV base installed cost (2022)
    23 = @(1492, 1492, 1492, 1492, 1492)   # PV base installed cost (2025)
    24 = @(1193, 1193, 1193, 1193, 1193)   # PV base installed cost (2030)
    25 = @(895, 895, 895, 895, 895)        # PV base installed cost (2035)
    26 = @(683, 683, 683, 683, 683)        # PV base installed cost (2050)

    28 = @(24, 24, 24, 24, 24)             # PV OpEx (2022)
    29 = @(21, 21, 21, 21, 21)             # PV OpEx (2025)
    30 = @(18, 18, 18, 18, 18)             # PV OpEx (2030)
    31 = @(15, 15, 15, 15, 15)             # PV OpEx (2035)
    32 = @(13, 13, 13, 13, 13)             # PV OpEx (2050)
}

foreach ($row in $rowData.Keys) {
    $values = $rowData[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}

# --- Update the selection to match the saved view state ---
[void]$ws.Range("F28:F32").Select()
